$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(3)
$tr2 = $sh.TextFrame2.TextRange
Write-Host "Text: [$($tr2.Text)] Length=$($tr2.Length)"
$sub = $tr2.Characters(1, $tr2.Length)
Write-Host "Sub Text: [$($sub.Text)]"
$sub.Text = "5/18/2018"
Write-Host "After: [$($tr2.Text)]"
